# Auto-generated Excel COM-interop script
# Applies numeric value updates to market-price / profit columns (H-N)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets,
# matching the 'scheduled runner' data refresh described in the commit.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 275.6154
$ws.Range("I2").Value = 185.27272
$ws.Range("J2").Value = 772.5
$ws.Range("K2").Value = 185.27272
$ws.Range("L2").Value = 772.5
$ws.Range("M2").Value = -72.27271999999999
$ws.Range("N2").Value = -998.5
$ws.Range("H11").Value = 13.75
$ws.Range("I11").Value = 13.75
$ws.Range("K11").Value = 13.75
$ws.Range("M11").Value = 126.25
$ws.Range("H16").Value = 1250
$ws.Range("J16").Value = 1250
$ws.Range("L16").Value = 1250
$ws.Range("N16").Value = -1710
$ws.Range("H20").Value = 2481
$ws.Range("I20").Value = 750
$ws.Range("J20").Value = 4212
$ws.Range("K20").Value = 750
$ws.Range("L20").Value = 4212
$ws.Range("M20").Value = -520
$ws.Range("N20").Value = -4672
$ws.Range("H35").Value = 2481
$ws.Range("I35").Value = 750
$ws.Range("J35").Value = 4212
$ws.Range("K35").Value = 750
$ws.Range("L35").Value = 4212
$ws.Range("M35").Value = -371
$ws.Range("N35").Value = -4970
$ws.Range("H68").Value = 48795
$ws.Range("J68").Value = 48795
$ws.Range("L68").Value = 48795
$ws.Range("N68").Value = -50293
$ws.Range("H71").Value = 48795
$ws.Range("J71").Value = 48795
$ws.Range("L71").Value = 146385
$ws.Range("N71").Value = -153873
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H131").Value = 720
$ws.Range("I131").Value = 720
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 2160
$ws.Range("L131").Value = 0
$ws.Range("M131").Value = 2880
$ws.Range("N131").ClearContents()
$ws.Range("H137").Value = 2707.7693
$ws.Range("I137").Value = 1033.5
$ws.Range("J137").Value = 4142.857
$ws.Range("K137").Value = 3100.5
$ws.Range("L137").Value = 12428.571
$ws.Range("M137").Value = -550.5
$ws.Range("N137").Value = -17528.571

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2473.125
$ws.Range("I61").Value = 2473.125
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2473.125
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -2261.125
$ws.Range("N61").ClearContents()
$ws.Range("H136").Value = 2473.125
$ws.Range("I136").Value = 2473.125
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 7419.375
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -4869.375
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H24").Value = 2800
$ws.Range("I24").Value = 3200
$ws.Range("J24").Value = 2000
$ws.Range("K24").Value = 3200
$ws.Range("L24").Value = 2000
$ws.Range("M24").Value = -2965
$ws.Range("N24").Value = -2470
$ws.Range("H88").Value = 18360.334
$ws.Range("J88").Value = 18360.334
$ws.Range("L88").Value = 18360.334
$ws.Range("N88").Value = -19172.334
$ws.Range("H91").Value = 18360.334
$ws.Range("J91").Value = 18360.334
$ws.Range("L91").Value = 18360.334
$ws.Range("N91").Value = -21168.334
$ws.Range("H92").Value = 50000
$ws.Range("J92").Value = 50000
$ws.Range("L92").Value = 50000
$ws.Range("N92").Value = -54992

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
$ws.Range("H92").Value = 37649.25
$ws.Range("J92").Value = 37649.25
$ws.Range("L92").Value = 37649.25
$ws.Range("N92").Value = -42641.25
$ws.Range("H96").Value = 6836.8
$ws.Range("J96").Value = 6836.8
$ws.Range("L96").Value = 6836.8
$ws.Range("N96").Value = -12328.8
$ws.Range("H99").Value = 2999.2856
$ws.Range("I99").Value = 2199.2
$ws.Range("J99").Value = 4999.5
$ws.Range("K99").Value = 2199.2
$ws.Range("L99").Value = 4999.5
$ws.Range("M99").Value = -701.1999999999998
$ws.Range("N99").Value = -7995.5
$ws.Range("H122").Value = 1245.5
$ws.Range("J122").Value = 1255.75
$ws.Range("L122").Value = 3767.25
$ws.Range("N122").Value = -8667.25
$ws.Range("H126").Value = 2999.2856
$ws.Range("I126").Value = 2199.2
$ws.Range("J126").Value = 4999.5
$ws.Range("K126").Value = 6597.599999999999
$ws.Range("L126").Value = 14998.5
$ws.Range("M126").Value = -4127.599999999999
$ws.Range("N126").Value = -19938.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 59.666668
$ws.Range("I2").Value = 116.333336
$ws.Range("J2").Value = 40.77778
$ws.Range("K2").Value = 698.000016
$ws.Range("L2").Value = 244.66668
$ws.Range("M2").Value = -585.000016
$ws.Range("N2").Value = -470.66668
$ws.Range("H12").Value = 205.57143
$ws.Range("I12").Value = 226
$ws.Range("J12").Value = 199.1875
$ws.Range("K12").Value = 678
$ws.Range("L12").Value = 597.5625
$ws.Range("M12").Value = -505
$ws.Range("N12").Value = -943.5625
$ws.Range("H40").Value = 57.6
$ws.Range("J40").Value = 111
$ws.Range("L40").Value = 444
$ws.Range("N40").Value = -582
$ws.Range("H86").Value = 100
$ws.Range("I86").Value = 100
$ws.Range("K86").Value = 300
$ws.Range("M86").Value = 886
$ws.Range("H89").Value = 100
$ws.Range("I89").Value = 100
$ws.Range("K89").Value = 900
$ws.Range("M89").Value = 5028
$ws.Range("H122").Value = 1060.2
$ws.Range("I122").Value = 1898
$ws.Range("J122").Value = 850.75
$ws.Range("K122").Value = 17082
$ws.Range("L122").Value = 7656.75
$ws.Range("M122").Value = -14632
$ws.Range("N122").Value = -12556.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 19186116
$ws.Range("I3").Value = 27962962
$ws.Range("J3").Value = 13334883
$ws.Range("K3").Value = 27962962
$ws.Range("L3").Value = 13334883
$ws.Range("M3").Value = -27962846
$ws.Range("N3").Value = -13335115
$ws.Range("H9").Value = 239.125
$ws.Range("I9").Value = 285.5
$ws.Range("J9").Value = 100
$ws.Range("K9").Value = 285.5
$ws.Range("L9").Value = 100
$ws.Range("M9").Value = -115.5
$ws.Range("N9").Value = -440
$ws.Range("H10").Value = 2916.6667
$ws.Range("J10").Value = 3375
$ws.Range("L10").Value = 3375
$ws.Range("N10").Value = -3713

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 740.93335
$ws.Range("I22").Value = 777.2222
$ws.Range("J22").Value = 686.5
$ws.Range("K22").Value = 777.2222
$ws.Range("L22").Value = 686.5
$ws.Range("M22").Value = -482.2222
$ws.Range("N22").Value = -1276.5
$ws.Range("H27").Value = 740.93335
$ws.Range("I27").Value = 777.2222
$ws.Range("J27").Value = 686.5
$ws.Range("K27").Value = 777.2222
$ws.Range("L27").Value = 686.5
$ws.Range("M27").Value = -670.2222
$ws.Range("N27").Value = -900.5
$ws.Range("H35").Value = 4426.2
$ws.Range("I35").Value = 4274
$ws.Range("K35").Value = 4274
$ws.Range("M35").Value = -3938
$ws.Range("H61").Value = 3339.8096
$ws.Range("I61").Value = 1438.4286
$ws.Range("K61").Value = 1438.4286
$ws.Range("M61").Value = -1236.4286
$ws.Range("H88").Value = 15000
$ws.Range("I88").Value = 15000
$ws.Range("K88").Value = 15000
$ws.Range("M88").Value = -14572
$ws.Range("H91").Value = 15000
$ws.Range("I91").Value = 15000
$ws.Range("K91").Value = 15000
$ws.Range("M91").Value = -13518
$ws.Range("H94").Value = 25000
$ws.Range("I94").Value = 25000
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 25000
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -24324
$ws.Range("N94").ClearContents()
$ws.Range("H113").Value = 3339.8096
$ws.Range("I113").Value = 1438.4286
$ws.Range("K113").Value = 1438.4286
$ws.Range("M113").Value = 731.5714

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").ClearContents()
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("H96").Value = 1750.2
$ws.Range("J96").Value = 1633
$ws.Range("L96").Value = 1633
$ws.Range("N96").Value = -4379
$ws.Range("H126").Value = 3175.8076
$ws.Range("I126").Value = 1351.2941
$ws.Range("J126").Value = 6622.1113
$ws.Range("K126").Value = 4053.8823
$ws.Range("L126").Value = 19866.3339
$ws.Range("M126").Value = -1583.8823
$ws.Range("N126").Value = -24806.3339
$ws.Range("H132").Value = 3528.4443
$ws.Range("I132").Value = 2712.5386
$ws.Range("K132").Value = 8137.6158
$ws.Range("M132").Value = -5607.6158

